$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: a date value (formatted as a date, numFmtId 14)
$dateCell = $ws.Range("B6")
$dateCell.Value = 46001
$dateCell.NumberFormat = "mm-dd-yy"

# Row 7: label + gradient (color-scale) cell
$ws.Range("A7").Value = "with gradient"

$gradientCell = $ws.Range("B7")
$gradientCell.Interior.Color = 16770303

$gradientCell.FormatConditions.AddColorScale(3) | Out-Null

# Move the active selection to the newly added cell, like Excel would after data entry
$ws.Range("B7").Select() | Out-Null
